# chore: update Sheets via scheduled runner
# Refreshes cached marketboard price snapshots (currentAveragePrice* / LevePrice* /
# LeveProfit* columns) for the affected leve rows across each job sheet.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1324.5
$ws.Range("I19").Value = 1473.75
$ws.Range("K19").Value = 1473.75
$ws.Range("M19").Value = -1298.75
$ws.Range("H53").Value = 364
$ws.Range("I53").Value = 344.5
$ws.Range("J53").Value = 422.5
$ws.Range("K53").Value = 344.5
$ws.Range("L53").Value = 422.5
$ws.Range("M53").Value = 292.5
$ws.Range("N53").Value = -1696.5
$ws.Range("H131").Value = 9499.333000000001
$ws.Range("I131").Value = 9578.200000000001
$ws.Range("K131").Value = 28734.6
$ws.Range("M131").Value = -23694.6
$ws.Range("H137").Value = 3799.8
$ws.Range("J137").Value = 3799.8
$ws.Range("L137").Value = 11399.4
$ws.Range("N137").Value = -16499.4
$ws.Range("H138").Value = 7253
$ws.Range("I138").Value = 2997.8
$ws.Range("J138").Value = 7878.7646
$ws.Range("K138").Value = 8993.400000000001
$ws.Range("L138").Value = 23636.2938
$ws.Range("M138").Value = -3853.400000000001
$ws.Range("N138").Value = -33916.2938

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2704.125
$ws.Range("I32").Value = 2662.3225
$ws.Range("K32").Value = 2662.3225
$ws.Range("M32").Value = -2375.3225
$ws.Range("H61").Value = 5003.778
$ws.Range("I61").Value = 5291
$ws.Range("K61").Value = 5291
$ws.Range("M61").Value = -5079
$ws.Range("H98").Value = 30355
$ws.Range("J98").Value = 30355
$ws.Range("L98").Value = 30355
$ws.Range("N98").Value = -36345
$ws.Range("H102").Value = 1300.5
$ws.Range("I102").Value = 1365.7
$ws.Range("K102").Value = 1365.7
$ws.Range("M102").Value = 256.3
$ws.Range("H132").Value = 3055.5
$ws.Range("I132").Value = 2176.3333
$ws.Range("J132").Value = 4374.25
$ws.Range("K132").Value = 6528.999899999999
$ws.Range("L132").Value = 13122.75
$ws.Range("M132").Value = -3998.999899999999
$ws.Range("N132").Value = -18182.75
$ws.Range("H136").Value = 5003.778
$ws.Range("I136").Value = 5291
$ws.Range("K136").Value = 15873
$ws.Range("M136").Value = -13323

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1450
$ws.Range("I99").Value = 1450
$ws.Range("K99").Value = 1450
$ws.Range("M99").Value = 48
$ws.Range("H134").Value = 3689.3572
$ws.Range("I134").Value = 3498.2
$ws.Range("J134").Value = 4167.25
$ws.Range("K134").Value = 10494.6
$ws.Range("L134").Value = 12501.75
$ws.Range("M134").Value = -7959.599999999999
$ws.Range("N134").Value = -17571.75
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 76.59999999999999
$ws.Range("I7").Value = 10
$ws.Range("J7").Value = 121
$ws.Range("K7").Value = 10
$ws.Range("L7").Value = 121
$ws.Range("M7").Value = 103
$ws.Range("N7").Value = -347
$ws.Range("H31").Value = 6125.9
$ws.Range("J31").Value = 9580.091
$ws.Range("L31").Value = 9580.091
$ws.Range("N31").Value = -10170.091
$ws.Range("H34").Value = 6125.9
$ws.Range("J34").Value = 9580.091
$ws.Range("L34").Value = 9580.091
$ws.Range("N34").Value = -9984.091
$ws.Range("H133").Value = 124847.5
$ws.Range("J133").Value = 124847.5
$ws.Range("L133").Value = 124847.5
$ws.Range("N133").Value = -129907.5
$ws.Range("H134").Value = 2210
$ws.Range("I134").Value = 2082.1428
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 6246.428400000001
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -3711.428400000001
$ws.Range("N134").Value = -17070

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1255.8572
$ws.Range("J68").Value = 1265.1666
$ws.Range("L68").Value = 3795.4998
$ws.Range("N68").Value = -5417.4998
$ws.Range("H71").Value = 1255.8572
$ws.Range("J71").Value = 1265.1666
$ws.Range("L71").Value = 11386.4994
$ws.Range("N71").Value = -19498.4994
$ws.Range("H103").Value = 362.66666
$ws.Range("I103").Value = 425
$ws.Range("J103").Value = 238
$ws.Range("K103").Value = 1275
$ws.Range("L103").Value = 714
$ws.Range("M103").Value = -396
$ws.Range("N103").Value = -2472
$ws.Range("H129").Value = 516.5
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 533
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 1599
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -11599
$ws.Range("H131").Value = 1004.8333
$ws.Range("I131").Value = 509.66666
$ws.Range("J131").Value = 1500
$ws.Range("K131").Value = 1528.99998
$ws.Range("L131").Value = 4500
$ws.Range("M131").Value = 3511.00002
$ws.Range("N131").Value = -14580

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 12999
$ws.Range("J23").Value = 12999
$ws.Range("L23").Value = 12999
$ws.Range("N23").Value = -13445
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H97").Value = 1067.8
$ws.Range("I97").Value = 575.2
$ws.Range("J97").Value = 2053
$ws.Range("K97").Value = 575.2
$ws.Range("L97").Value = 2053
$ws.Range("M97").Value = -79.20000000000005
$ws.Range("N97").Value = -3045
$ws.Range("H126").Value = 2148.6
$ws.Range("I126").Value = 2299
$ws.Range("J126").Value = 1998.2
$ws.Range("K126").Value = 6897
$ws.Range("L126").Value = 5994.6
$ws.Range("M126").Value = -4427
$ws.Range("N126").Value = -10934.6
$ws.Range("H132").Value = 3309
$ws.Range("I132").Value = 2760.4285
$ws.Range("J132").Value = 4749
$ws.Range("K132").Value = 8281.2855
$ws.Range("L132").Value = 14247
$ws.Range("M132").Value = -5751.2855
$ws.Range("N132").Value = -19307

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2892.1875
$ws.Range("I40").Value = 2791.6428
$ws.Range("K40").Value = 2791.6428
$ws.Range("M40").Value = -2655.6428
$ws.Range("H122").Value = 1574.5
$ws.Range("I122").Value = 1574.5
$ws.Range("K122").Value = 4723.5
$ws.Range("M122").Value = -2273.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1332.6666
$ws.Range("I122").Value = 1249.3889
$ws.Range("J122").Value = 1832.3334
$ws.Range("K122").Value = 3748.1667
$ws.Range("L122").Value = 5497.0002
$ws.Range("M122").Value = -1298.1667
$ws.Range("N122").Value = -10397.0002
$ws.Range("H136").Value = 5250.9653
$ws.Range("I136").Value = 6170.304
$ws.Range("K136").Value = 18510.912
$ws.Range("M136").Value = -15960.912
